# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates the "K" column (column G) values for each saved game row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row -> new K value
$kValues = @{
    2  = 0
    3  = 3
    4  = 2
    5  = 2
    6  = 0
    7  = 2
    8  = 4
    9  = 0
    10 = 3
    11 = 2
    12 = 2
    13 = 4
    14 = 2
    15 = 2
    16 = 2
    17 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
